# Updates crypto price (D) and 1h-volume-change (E) columns to match latest
# scrape, per the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.961.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.381.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.19%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.565'
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.510'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.62'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0789'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.64%  '

$ws.Range("E13").Value = '  -0.40%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.745.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.19%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.390.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.828'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.50%  '

$ws.Range("E17").Value = '  -2.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.867.63'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0953'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.73%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '39.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.98%  '

$ws.Range("E28").Value = '  -3.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.80'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +18.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.93%  '

$ws.Range("E32").Value = '  +7.09%  '

$ws.Range("E33").Value = '  -4.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '146.36'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0774'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.112'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("E37").Value = '  +5.95%  '

$ws.Range("E38").Value = '  -2.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.70%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0300'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.25'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.82%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.933.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.54%  '

$ws.Range("E44").Value = '  +0.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.28%  '

$ws.Range("E46").Value = '  -10.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.16%  '

$ws.Range("E48").Value = '  -5.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '98.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.615.02'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '68.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.36%  '
